$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "Agen" column (column I) - removes the supplier field
$ws.Columns("I").Delete()

# Update selection to match the post-edit state
$ws.Range("I1").Select()
